# Auto-generated files on 2025-12-28
# Refresh the daily "Hot Stock Top 20" rankings on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (财联社) - new ranking, rows 2-21
$colA = @(
    "航天发展",
    "神剑股份",
    "中国卫星",
    "中超控股",
    "锋龙股份",
    "通宇通讯",
    "天际股份",
    "中国卫通",
    "航天电子",
    "再升科技",
    "泰尔股份",
    "金风科技",
    "江西铜业",
    "海南发展",
    "北斗星通",
    "平潭发展",
    "天奇股份",
    "胜通能源",
    "白银有色",
    "安通控股"
)

# Column B (东方财富) - new ranking, rows 2-21
$colB = @(
    "中超控股",
    "江西铜业",
    "航天发展",
    "神剑股份",
    "锋龙股份",
    "中国卫星",
    "泰尔股份",
    "白银有色",
    "通宇通讯",
    "中国卫通",
    "海南发展",
    "金风科技",
    "航天电子",
    "天际股份",
    "北斗星通",
    "再升科技",
    "福龙马",
    "永鼎股份",
    "安通控股",
    "天奇股份"
)

# Column C (同花顺) - new ranking, rows 2-21
$colC = @(
    "航天发展",
    "中超控股",
    "神剑股份",
    "平潭发展",
    "锋龙股份",
    "中国卫星",
    "通宇通讯",
    "再升科技",
    "胜通能源",
    "天际股份",
    "泰尔股份",
    "海南发展",
    "雪人集团",
    "浙江世宝",
    "东百集团",
    "安通控股",
    "航天电子",
    "嘉美包装",
    "金风科技",
    "中国卫通"
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
    $ws.Cells.Item($row, 3).Value = $colC[$i]
}
